{"js": "// Replace the inline picture (the \"dracula\" image) in the first paragraph\n// with real text: \"Test document \". Also drop the now-unused \"Balloon Text\"\n// / \"Balloon Text Char\" styles from the style sheet.\n\n// 1) Swap the inline picture for plain text in-place.\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length > 0) {\n  const picture = pictures.items[0];\n  const pictureRange = picture.getRange();\n  pictureRange.insertText(\"Test document \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the unused \"Balloon Text\" / \"Balloon Text Char\" styles.\n//    Delete the linked character style first, then its base paragraph\n//    style (deleting in the other order can leave the collection in an\n//    inconsistent state).\nconst charStyle = context.document.getStyles().getByNameOrNullObject(\"Balloon Text Char\");\ncharStyle.load(\"nameLocal\");\nawait context.sync();\nif (!charStyle.isNullObject) {\n  charStyle.delete();\n  await context.sync();\n}\n\nconst paragraphStyle = context.document.getStyles().getByNameOrNullObject(\"Balloon Text\");\nparagraphStyle.load(\"nameLocal\");\nawait context.sync();\nif (!paragraphStyle.isNullObject) {\n  paragraphStyle.delete();\n  await context.sync();\n}\n", "ps1": "# Replace the inline picture (the \"dracula\" image) in the first paragraph\n# with real text: \"Test document \". Also drop the now-unused \"Balloon Text\"\n# / \"Balloon Text Char\" styles from the style sheet.\n\n$d = $word.ActiveDocument\n\n# 1) Swap the inline picture for plain text in-place.\nif ($d.InlineShapes.Count -gt 0) {\n    $shape = $d.InlineShapes(1)\n    $rng = $shape.Range\n    $rng.Text = \"Test document \"\n}\n\n# 2) Remove the unused \"Balloon Text\" / \"Balloon Text Char\" styles.\n#    Delete the linked character style first, then its base paragraph\n#    style (deleting in the other order can leave the collection in an\n#    inconsistent state). Guard with an existence check so the script is\n#    safe to run even if the styles were already removed.\n$styleNames = @()\nforeach ($s in $d.Styles) {\n    $styleNames += $s.NameLocal\n}\n\nif ($styleNames -contains \"Balloon Text Char\") {\n    $d.Styles(\"Balloon Text Char\").Delete()\n}\nif ($styleNames -contains \"Balloon Text\") {\n    $d.Styles(\"Balloon Text\").Delete()\n}\n"}
